$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '60.979.61'
$ws.Range("E2").Value = '  -3.25%  '
Set-TextValue $ws.Range("D3") '3.001.35'
$ws.Range("E3").Value = '  -2.28%  '
$ws.Range("E4").Value = '  +0.14%  '
Set-TextValue $ws.Range("D5") '531.93'
$ws.Range("E5").Value = '  -0.92%  '
Set-TextValue $ws.Range("D6") '133.25'
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("E7").Value = '  -0.03%  '
Set-TextValue $ws.Range("D8") '3.001.00'
$ws.Range("E8").Value = '  -2.04%  '
Set-TextValue $ws.Range("D9") '0.497'
$ws.Range("E9").Value = '  +0.87%  '
Set-TextValue $ws.Range("D10") '0.148'
$ws.Range("E10").Value = '  -3.90%  '
Set-TextValue $ws.Range("D11") '6.06'
$ws.Range("E11").Value = '  -1.46%  '
Set-TextValue $ws.Range("D12") '0.445'
$ws.Range("E12").Value = '  -1.42%  '
Set-TextValue $ws.Range("D13") '0.0000220'
$ws.Range("E13").Value = '  -1.52%  '
Set-TextValue $ws.Range("D14") '34.15'
$ws.Range("E14").Value = '  -0.24%  '
Set-TextValue $ws.Range("D15") '3.501.22'
$ws.Range("E15").Value = '  -1.69%  '
$ws.Range("E16").Value = '  -0.40%  '
Set-TextValue $ws.Range("D17") '61.148.23'
$ws.Range("E17").Value = '  -2.91%  '
Set-TextValue $ws.Range("D18") '3.011.58'
$ws.Range("E18").Value = '  -1.86%  '
Set-TextValue $ws.Range("D19") '6.58'
$ws.Range("E19").Value = '  -0.51%  '
Set-TextValue $ws.Range("D20") '462.34'
$ws.Range("E20").Value = '  -4.00%  '
Set-TextValue $ws.Range("D21") '13.22'
$ws.Range("E21").Value = '  -0.70%  '
Set-TextValue $ws.Range("D22") '0.673'
$ws.Range("E22").Value = '  -2.86%  '
Set-TextValue $ws.Range("D23") '6.92'
$ws.Range("E23").Value = '  -2.43%  '
Set-TextValue $ws.Range("D24") '79.19'
$ws.Range("E24").Value = '  +0.34%  '
Set-TextValue $ws.Range("D25") '12.00'
$ws.Range("E25").Value = '  -0.73%  '
$ws.Range("E26").Value = '  +0.10%  '
Set-TextValue $ws.Range("D27") '2.67'
$ws.Range("E27").Value = '  -1.11%  '
$ws.Range("E28").Value = '  -2.86%  '
Set-TextValue $ws.Range("D29") '1.00'
$ws.Range("E29").Value = '  +0.43%  '
Set-TextValue $ws.Range("D30") '1.89'
$ws.Range("E30").Value = '  +0.35%  '
Set-TextValue $ws.Range("D31") '25.45'
$ws.Range("E31").Value = '  -1.92%  '
$ws.Range("E32").Value = '  +2.88%  '
$ws.Range("B33").Value = 'OKB'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range("D33") '55.40'
$ws.Range("E33").Value = '  -2.67%  '
$ws.Range("B34").Value = 'NEARProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range("D34") '5.45'
$ws.Range("E34").Value = '  +2.07%  '
Set-TextValue $ws.Range("D35") '2.28'
$ws.Range("E35").Value = '  -3.37%  '
Set-TextValue $ws.Range("D36") '5.87'
$ws.Range("E36").Value = '  -2.34%  '
Set-TextValue $ws.Range("D37") '456.60'
$ws.Range("E37").Value = '  -4.43%  '
Set-TextValue $ws.Range("D38") '3.214.64'
$ws.Range("E38").Value = '  +4.00%  '
Set-TextValue $ws.Range("D39") '0.0783'
$ws.Range("E39").Value = '  -1.44%  '
Set-TextValue $ws.Range("D40") '0.0382'
$ws.Range("E40").Value = '  -2.76%  '
$ws.Range("E41").Value = '  +2.11%  '
Set-TextValue $ws.Range("D42") '8.13'
$ws.Range("E42").Value = '  +0.80%  '
Set-TextValue $ws.Range("D43") '27.49'
$ws.Range("E43").Value = '  +12.95%  '
Set-TextValue $ws.Range("D44") '2.46'
$ws.Range("E44").Value = '  -6.02%  '
$ws.Range("E45").Value = '  +0.10%  '
Set-TextValue $ws.Range("D46") '0.245'
$ws.Range("E46").Value = '  -2.62%  '
$ws.Range("E47").Value = '  -0.68%  '
Set-TextValue $ws.Range("D48") '118.81'
$ws.Range("E48").Value = '  -2.22%  '
$ws.Range("E49").Value = '  +0.26%  '
Set-TextValue $ws.Range("D50") '0.0₃0492'
$ws.Range("E50").Value = '  -8.77%  '
$ws.Range("E51").Value = '  +7.66%  '
